$d = $word.ActiveDocument

# Locate the empty paragraph immediately following the HUIT-019 user
# story block (identified by a unique, accent-free anchor phrase near
# its end) -- this is the paragraph that gets replaced by the new
# HUIT-020 story plus a fresh trailing empty paragraph.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Evitar confusiones", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the HUIT-019 anchor text"
}

$targetIndex = -1
$bestStart = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $s = $p.Range.Start
    if ($s -ge $rng.End) {
        if ($targetIndex -eq -1 -or $s -lt $bestStart) {
            $targetIndex = $i
            $bestStart = $s
        }
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the target empty paragraph after HUIT-019"
}

$targetRange = $d.Paragraphs.Item($targetIndex).Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr><w:t>HUIT-002 Clasificación de publicidad por tipo</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr><w:br/><w:t>Como: Administrador</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr><w:br/><w:t>Quiero: Clasificar la publicidad por tipo, como nuevo lanzamiento, descuento, etc.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr><w:br/><w:t>Para: Llevar un mejor control y facilitar la búsqueda de la publicidad.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-GT"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($xml) | Out-Null
